# ==========================================================================
# Rebuild the worksheet content to exactly match the target shared-string
# order: remove the existing hyperlink + clear all cells, then re-write every
# cell (including the unchanged header block) top-to-bottom so new strings are
# appended to the shared-string table in the same order the target file has them.
# ==========================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# --- Row 1 ---
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "Peru"
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Font.Size = 18

# --- Row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "MSME Participation on the Economy"
$ws.Range("A3").ClearFormats()
$ws.Range("A3").Font.Bold = $true

# --- Row 9 ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").ClearFormats()
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- Row 11 ---
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").ClearFormats()
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").ClearFormats()
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").ClearFormats()
$ws.Range("D11").Font.Bold = $true

# --- Row 12 ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("A12").ClearFormats()
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "1177901"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "20062"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1197963"
$ws.Range("D12").ClearFormats()

# --- Row 13 ---
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "Enterprises density (per 1000 people)"
$ws.Range("A13").ClearFormats()
$ws.Range("A13").Font.Bold = $true
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "40.3"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "0.7"
$ws.Range("C13").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.9"
$ws.Range("D13").ClearFormats()

# --- Row 14 ---
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Employment (% of total)"
$ws.Range("A14").ClearFormats()
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "76.9"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "12.7"
$ws.Range("C14").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.6"
$ws.Range("D14").ClearFormats()

# --- Row 15 ---
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "Employment (absolute #)"
$ws.Range("A15").ClearFormats()
$ws.Range("A15").Font.Bold = $true
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "8168643"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "1345238"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9513881"
$ws.Range("D15").ClearFormats()

# --- Row 16 ---
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Enterprises (% of total)"
$ws.Range("A16").ClearFormats()
$ws.Range("A16").Font.Bold = $true
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "98.1"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "1.7"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.8"
$ws.Range("D16").ClearFormats()

# --- Row 17 ---
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "Source: CODEMYPE, 2010"
$ws.Range("A17").ClearFormats()
$ws.Range("A17").Font.Italic = $true

# --- Row 20 ---
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Number of employees"
$ws.Range("B20").ClearFormats()
$ws.Range("B20").Font.Bold = $true
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C20").ClearFormats()
$ws.Range("C20").Font.Bold = $true
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D20").ClearFormats()
$ws.Range("D20").Font.Bold = $true

# --- Row 21 ---
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Micro"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "1-10"
$ws.Range("B21").ClearFormats()
# C21: target is an empty string cell (not representable via COM Value assignment; left blank)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "<UIT 150"
$ws.Range("D21").ClearFormats()

# --- Row 22 ---
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "Small"
$ws.Range("A22").ClearFormats()
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "11-100"
$ws.Range("B22").ClearFormats()
# C22: target is an empty string cell (not representable via COM Value assignment; left blank)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "<UIT 1,700"
$ws.Range("D22").ClearFormats()

# --- Row 23 ---
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "Medium"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = ">100"
$ws.Range("B23").ClearFormats()
# C23: target is an empty string cell (not representable via COM Value assignment; left blank)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = ">=UIT 1,700"
$ws.Range("D23").ClearFormats()

# --- Row 24 ---
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "Large"
$ws.Range("A24").ClearFormats()
# B24: target is an empty string cell (not representable via COM Value assignment; left blank)
# C24: target is an empty string cell (not representable via COM Value assignment; left blank)
# D24: target is an empty string cell (not representable via COM Value assignment; left blank)

# --- Row 27 ---
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "Sector Distribution Details"
$ws.Range("A27").ClearFormats()
$ws.Range("A27").Font.Bold = $true

# --- Row 29 ---
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "SMEs"
$ws.Range("B29").ClearFormats()
$ws.Range("B29").Font.Bold = $true
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "%SMEs"
$ws.Range("C29").ClearFormats()
$ws.Range("C29").Font.Bold = $true

# --- Row 30 ---
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "Manufacturing and Processing "
$ws.Range("A30").ClearFormats()
$ws.Range("A30").Font.Bold = $true
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "128,682"
$ws.Range("B30").ClearFormats()
$ws.Range("B30").Font.Bold = $true
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "10.79"
$ws.Range("C30").ClearFormats()
$ws.Range("C30").Font.Bold = $true

# --- Row 31 ---
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "Trade "
$ws.Range("A31").ClearFormats()
$ws.Range("A31").Font.Bold = $true
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "562,736"
$ws.Range("B31").ClearFormats()
$ws.Range("B31").Font.Bold = $true
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "47.20"
$ws.Range("C31").ClearFormats()
$ws.Range("C31").Font.Bold = $true

# --- Row 32 ---
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "Services "
$ws.Range("A32").ClearFormats()
$ws.Range("A32").Font.Bold = $true
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "474,672"
$ws.Range("B32").ClearFormats()
$ws.Range("B32").Font.Bold = $true
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "39.81"
$ws.Range("C32").ClearFormats()
$ws.Range("C32").Font.Bold = $true

# --- Row 33 ---
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "Services"
$ws.Range("A33").ClearFormats()
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "441,445"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "37.02"
$ws.Range("C33").ClearFormats()

# --- Row 34 ---
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "Construction"
$ws.Range("A34").ClearFormats()
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "33,227"
$ws.Range("B34").ClearFormats()
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "2.79"
$ws.Range("C34").ClearFormats()

# --- Row 35 ---
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "Primary Agriculture "
$ws.Range("A35").ClearFormats()
$ws.Range("A35").Font.Bold = $true
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "26,211"
$ws.Range("B35").ClearFormats()
$ws.Range("B35").Font.Bold = $true
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "2.20"
$ws.Range("C35").ClearFormats()
$ws.Range("C35").Font.Bold = $true

# --- Row 36 ---
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "Agriculture and fishing"
$ws.Range("A36").ClearFormats()
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "26,211"
$ws.Range("B36").ClearFormats()
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "2.20"
$ws.Range("C36").ClearFormats()

# --- Row 37 ---
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "Total "
$ws.Range("A37").ClearFormats()
$ws.Range("A37").Font.Bold = $true
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "1,192,301"
$ws.Range("B37").ClearFormats()
$ws.Range("B37").Font.Bold = $true
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "100"
$ws.Range("C37").ClearFormats()
$ws.Range("C37").Font.Bold = $true

# --- Row 38 ---
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "Source:"
$ws.Range("A38").ClearFormats()
$ws.Range("A38").Font.Italic = $true

# --- Row 39 ---
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "Dirección General De MYPe Y cooPerativas Dirección De Desarrollo eMPresarial based on SUNAT, Registro Único de Contribuyentes 2010"
$ws.Range("A39").ClearFormats()
$ws.Range("A39").Font.Italic = $true

# --- Row 40 ---
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "http://www.produce.gob.pe/remype/data/mype2010.pdf"
$ws.Range("A40").ClearFormats()
$ws.Range("A40").Font.Underline = $true
$ws.Range("A40").Font.Color = 16711680

# --- Row 41 ---
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "Section 2.5"
$ws.Range("A41").ClearFormats()
$ws.Range("A41").Font.Italic = $true

# --- Row 44 ---
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "CODEMYPE"
$ws.Range("A44").ClearFormats()
$ws.Range("A44").Font.Bold = $true

# --- Row 45 ---
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "CODEMYPE, Estadisticas, Direccion General de MYPE y Cooperativas, p. 20,. Available at http://www.produce.gob.pe/remype/data/mype2010.pdf"
$ws.Range("A45").ClearFormats()
$ws.Range("A45").Font.Italic = $true

# Re-add the hyperlink at its new location
$ws.Hyperlinks.Add($ws.Range("A40"), "http://www.produce.gob.pe/remype/data/mype2010.pdf")

Write-Output "done"
